$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 516.6786
$ws.Range("I19").Value = 413.53845
$ws.Range("J19").Value = 606.06665
$ws.Range("K19").Value = 413.53845
$ws.Range("L19").Value = 606.06665
$ws.Range("M19").Value = -238.53845
$ws.Range("N19").Value = -956.06665
$ws.Range("H28").Value = 5052.5
$ws.Range("I28").Value = 5052.5
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 5052.5
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -4567.5
$ws.Range("N28").ClearContents()
$ws.Range("H58").Value = 1824.2174
$ws.Range("I58").Value = 328.66666
$ws.Range("J58").Value = 2785.6428
$ws.Range("K58").Value = 985.9999799999999
$ws.Range("L58").Value = 8356.928400000001
$ws.Range("M58").Value = -835.9999799999999
$ws.Range("N58").Value = -8656.928400000001
$ws.Range("H129").Value = 1236.9773
$ws.Range("I129").Value = 700
$ws.Range("J129").Value = 1249.4651
$ws.Range("K129").Value = 2100
$ws.Range("L129").Value = 3748.3953
$ws.Range("M129").Value = 2900
$ws.Range("N129").Value = -13748.3953
$ws.Range("H138").Value = 2350.2876
$ws.Range("I138").Value = 1259.9512
$ws.Range("J138").Value = 3496.5386
$ws.Range("K138").Value = 3779.8536
$ws.Range("L138").Value = 10489.6158
$ws.Range("M138").Value = 1360.1464
$ws.Range("N138").Value = -20769.6158
$ws.Range("H141").Value = 4421.2163
$ws.Range("I141").Value = 1508.9062
$ws.Range("J141").Value = 23060
$ws.Range("K141").Value = 4526.7186
$ws.Range("L141").Value = 69180
$ws.Range("M141").Value = 653.2813999999998
$ws.Range("N141").Value = -79540

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3850
$ws.Range("I88").Value = 2800
$ws.Range("J88").Value = 7000
$ws.Range("K88").Value = 2800
$ws.Range("L88").Value = 7000
$ws.Range("M88").Value = -2394
$ws.Range("N88").Value = -7812
$ws.Range("H91").Value = 3850
$ws.Range("I91").Value = 2800
$ws.Range("J91").Value = 7000
$ws.Range("K91").Value = 2800
$ws.Range("L91").Value = 7000
$ws.Range("M91").Value = -1396
$ws.Range("N91").Value = -9808

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 19190.291
$ws.Range("I107").Value = 66037.375
$ws.Range("K107").Value = 66037.375
$ws.Range("M107").Value = -64117.375

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1426840.4
$ws.Range("I58").Value = 2471409.2
$ws.Range("J58").Value = 2428.5454
$ws.Range("K58").Value = 2471409.2
$ws.Range("L58").Value = 2428.5454
$ws.Range("M58").Value = -2471206.2
$ws.Range("N58").Value = -2834.5454
$ws.Range("H107").Value = 1005.5
$ws.Range("I107").Value = 1005.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1005.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 914.5
$ws.Range("N107").ClearContents()
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H121").Value = 30000
$ws.Range("J121").Value = 30000
$ws.Range("L121").Value = 30000
$ws.Range("N121").Value = -32620
$ws.Range("H136").Value = 1426840.4
$ws.Range("I136").Value = 2471409.2
$ws.Range("J136").Value = 2428.5454
$ws.Range("K136").Value = 7414227.600000001
$ws.Range("L136").Value = 7285.6362
$ws.Range("M136").Value = -7411677.600000001
$ws.Range("N136").Value = -12385.6362

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2314.182
$ws.Range("I5").Value = 2932.75
$ws.Range("J5").Value = 664.6667
$ws.Range("K5").Value = 8798.25
$ws.Range("L5").Value = 1994.0001
$ws.Range("M5").Value = -8686.25
$ws.Range("N5").Value = -2218.0001
$ws.Range("H122").Value = 874.619
$ws.Range("I122").Value = 553.1111
$ws.Range("J122").Value = 1115.75
$ws.Range("K122").Value = 4977.9999
$ws.Range("L122").Value = 10041.75
$ws.Range("M122").Value = -2527.9999
$ws.Range("N122").Value = -14941.75
$ws.Range("H131").Value = 15153121
$ws.Range("J131").Value = 15386239
$ws.Range("L131").Value = 46158717
$ws.Range("N131").Value = -46168797
$ws.Range("H135").Value = 2314.182
$ws.Range("I135").Value = 2932.75
$ws.Range("J135").Value = 664.6667
$ws.Range("K135").Value = 26394.75
$ws.Range("L135").Value = 5982.0003
$ws.Range("M135").Value = -23859.75
$ws.Range("N135").Value = -11052.0003

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1386.1613
$ws.Range("I113").Value = 1109.8462
$ws.Range("J113").Value = 1585.7222
$ws.Range("K113").Value = 1109.8462
$ws.Range("L113").Value = 1585.7222
$ws.Range("M113").Value = 1060.1538
$ws.Range("N113").Value = -5925.7222
$ws.Range("H124").Value = 49780
$ws.Range("J124").Value = 49780
$ws.Range("L124").Value = 49780
$ws.Range("N124").Value = -59600

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3849.9
$ws.Range("J68").Value = 3499.875
$ws.Range("L68").Value = 3499.875
$ws.Range("N68").Value = -4997.875
$ws.Range("H71").Value = 3849.9
$ws.Range("J71").Value = 3499.875
$ws.Range("L71").Value = 17499.375
$ws.Range("N71").Value = -24987.375
$ws.Range("H101").Value = 173453.33
$ws.Range("J101").Value = 173453.33
$ws.Range("L101").Value = 173453.33
$ws.Range("N101").Value = -179943.33
$ws.Range("H121").Value = 75166.164
$ws.Range("J121").Value = 75166.164
$ws.Range("L121").Value = 75166.164
$ws.Range("N121").Value = -78660.164
$ws.Range("H122").Value = 18186164
$ws.Range("I122").Value = 4628.5713
$ws.Range("K122").Value = 13885.7139
$ws.Range("M122").Value = -11435.7139
$ws.Range("H132").Value = 3598.3684
$ws.Range("I132").Value = 2712.1428
$ws.Range("J132").Value = 4115.3335
$ws.Range("K132").Value = 8136.428400000001
$ws.Range("L132").Value = 12346.0005
$ws.Range("M132").Value = -5606.428400000001
$ws.Range("N132").Value = -17406.0005

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 22285.715
$ws.Range("J103").Value = 22285.715
$ws.Range("L103").Value = 22285.715
$ws.Range("N103").Value = -24629.715
$ws.Range("H112").Value = 500000
$ws.Range("J112").Value = 500000
$ws.Range("L112").Value = 500000
$ws.Range("N112").Value = -502954
$ws.Range("H113").Value = 737.5
$ws.Range("I113").Value = 721.8182
$ws.Range("J113").Value = 795
$ws.Range("K113").Value = 2165.4546
$ws.Range("L113").Value = 2385
$ws.Range("M113").Value = 4.545399999999972
$ws.Range("N113").Value = -6725
$ws.Range("H121").Value = 32481.889
$ws.Range("J121").Value = 32481.889
$ws.Range("L121").Value = 32481.889
$ws.Range("N121").Value = -35975.889
$ws.Range("H123").Value = 22238.172
$ws.Range("J123").Value = 22238.172
$ws.Range("L123").Value = 22238.172
$ws.Range("N123").Value = -32038.172
